# Projection info is added
# Calculate the generated shapefile file space coordinate information supplement.
#
# Updates the computed indices for rows 2-3 in columns G (VF) and H (Iat),
# and narrows column G's width to match the other data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (VF) is 12.7109375 wide in the source; narrow it to match the
# other feature columns (11.7109375). The COM layer only accepts
# ColumnWidth in 1/6-pt increments, so 10.8 is the closest input that
# lands on the nearest achievable stored width.
$ws.Columns.Item(7).ColumnWidth = 10.8

# Recomputed coordinate values for the VF (G) and Iat (H) columns.
$ws.Range("G2").Value = 2.955188512802124
$ws.Range("H2").Value = 2.1666667461395264
$ws.Range("G3").Value = 3.530733585357666
$ws.Range("H3").Value = 2.5
